# Updated, more de-identified profiles
# De-identify the "Visited U.S.?" column (D) by collapsing detailed
# location write-ins down to a plain "Yes", matching the author's commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value  = "Yes"
$ws.Range("D6").Value  = "Yes"
$ws.Range("D8").Value  = "Yes"
$ws.Range("D10").Value = "Yes"
$ws.Range("D11").Value = "Yes"
$ws.Range("D12").Value = "Yes"
$ws.Range("D14").Value = "Yes"

# Reflect the saved view/selection state (scrolled so row 3 is at top,
# active cell D13 selected).
$ws.Range("D13").Select()
$excel.ActiveWindow.ScrollRow = 3
